# Update Name of Algo
# Apply updated KNN imputation result values to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 5.654
$ws.Range("C4").Value = -12.74

$ws.Range("C5").Value = -12.927

$ws.Range("B6").Value = 7.295999999999999

$ws.Range("B7").Value = 6.611999999999999

$ws.Range("C8").Value = -12.8

$ws.Range("B16").Value = 5.801
$ws.Range("C16").Value = -12.523

$ws.Range("B20").Value = 6.462000000000001

$ws.Range("C22").Value = -12.78
